$wb = $excel.ActiveWorkbook

# --- Sheet1: LP1912 ---
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = "Última actualización: 17:13:30"
$ws1.Range("A3").Value = "Total filas: 44"

$ws1.Cells.Item(18, 1).Value = "17:13:30"
$ws1.Cells.Item(18, 2).Value = "17:16"
$ws1.Cells.Item(18, 3).Value = "10_OLMOS"
$ws1.Cells.Item(18, 4).Value = 3
$ws1.Cells.Item(18, 5).Value = "LP1912"

$ws1.Cells.Item(19, 1).Value = "16:50:41"
$ws1.Cells.Item(19, 2).Value = "17:17"
$ws1.Cells.Item(19, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(19, 4).Value = 27
$ws1.Cells.Item(19, 5).Value = "LP1912"

$ws1.Cells.Item(20, 1).Value = "16:46:42"
$ws1.Cells.Item(20, 2).Value = "17:17"
$ws1.Cells.Item(20, 3).Value = "17_ROMERO"
$ws1.Cells.Item(20, 4).Value = 31
$ws1.Cells.Item(20, 5).Value = "LP1912"

$ws1.Cells.Item(21, 1).Value = "16:52:37"
$ws1.Cells.Item(21, 2).Value = "17:20"
$ws1.Cells.Item(21, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(21, 4).Value = 28
$ws1.Cells.Item(21, 5).Value = "LP1912"

$ws1.Cells.Item(22, 1).Value = "16:46:42"
$ws1.Cells.Item(22, 2).Value = "17:23"
$ws1.Cells.Item(22, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(22, 4).Value = 37
$ws1.Cells.Item(22, 5).Value = "LP1912"

$ws1.Cells.Item(23, 1).Value = "16:46:42"
$ws1.Cells.Item(23, 2).Value = "17:24"
$ws1.Cells.Item(23, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(23, 4).Value = 38
$ws1.Cells.Item(23, 5).Value = "LP1912"

$ws1.Cells.Item(24, 1).Value = "17:13:30"
$ws1.Cells.Item(24, 2).Value = "17:27"
$ws1.Cells.Item(24, 3).Value = "15_ABASTO"
$ws1.Cells.Item(24, 4).Value = 14
$ws1.Cells.Item(24, 5).Value = "LP1912"

$ws1.Cells.Item(25, 1).Value = "17:13:30"
$ws1.Cells.Item(25, 2).Value = "17:33"
$ws1.Cells.Item(25, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(25, 4).Value = 20
$ws1.Cells.Item(25, 5).Value = "LP1912"

$ws1.Cells.Item(26, 1).Value = "16:50:41"
$ws1.Cells.Item(26, 2).Value = "17:34"
$ws1.Cells.Item(26, 3).Value = "10_OLMOS"
$ws1.Cells.Item(26, 4).Value = 44
$ws1.Cells.Item(26, 5).Value = "LP1912"

$ws1.Cells.Item(27, 1).Value = "16:46:42"
$ws1.Cells.Item(27, 2).Value = "17:35"
$ws1.Cells.Item(27, 3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(27, 4).Value = 49
$ws1.Cells.Item(27, 5).Value = "LP1912"

$ws1.Cells.Item(28, 1).Value = "16:52:37"
$ws1.Cells.Item(28, 2).Value = "17:36"
$ws1.Cells.Item(28, 3).Value = "27_EL RETIRO"
$ws1.Cells.Item(28, 4).Value = 44
$ws1.Cells.Item(28, 5).Value = "LP1912"

$ws1.Cells.Item(29, 1).Value = "17:13:30"
$ws1.Cells.Item(29, 2).Value = "17:37"
$ws1.Cells.Item(29, 3).Value = "27_EL RETIRO"
$ws1.Cells.Item(29, 4).Value = 24
$ws1.Cells.Item(29, 5).Value = "LP1912"

$ws1.Cells.Item(30, 1).Value = "16:46:42"
$ws1.Cells.Item(30, 2).Value = "17:38"
$ws1.Cells.Item(30, 3).Value = "17X38_ROMERO"
$ws1.Cells.Item(30, 4).Value = 52
$ws1.Cells.Item(30, 5).Value = "LP1912"

$ws1.Cells.Item(31, 1).Value = "17:13:30"
$ws1.Cells.Item(31, 2).Value = "17:41"
$ws1.Cells.Item(31, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(31, 4).Value = 28
$ws1.Cells.Item(31, 5).Value = "LP1912"

$ws1.Cells.Item(32, 1).Value = "16:46:42"
$ws1.Cells.Item(32, 2).Value = "17:44"
$ws1.Cells.Item(32, 3).Value = "215B_EL PATO"
$ws1.Cells.Item(32, 4).Value = 58
$ws1.Cells.Item(32, 5).Value = "LP1912"

$ws1.Cells.Item(33, 1).Value = "16:50:41"
$ws1.Cells.Item(33, 2).Value = "17:47"
$ws1.Cells.Item(33, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(33, 4).Value = 57
$ws1.Cells.Item(33, 5).Value = "LP1912"

$ws1.Cells.Item(34, 1).Value = "16:46:42"
$ws1.Cells.Item(34, 2).Value = "17:48"
$ws1.Cells.Item(34, 3).Value = "27_EL RETIRO"
$ws1.Cells.Item(34, 4).Value = 62
$ws1.Cells.Item(34, 5).Value = "LP1912"

$ws1.Cells.Item(35, 1).Value = "16:50:41"
$ws1.Cells.Item(35, 2).Value = "17:49"
$ws1.Cells.Item(35, 3).Value = "27_EL RETIRO"
$ws1.Cells.Item(35, 4).Value = 59
$ws1.Cells.Item(35, 5).Value = "LP1912"

$ws1.Cells.Item(36, 1).Value = "16:46:42"
$ws1.Cells.Item(36, 2).Value = "17:50"
$ws1.Cells.Item(36, 3).Value = "215_EL PELIGRO"
$ws1.Cells.Item(36, 4).Value = 64
$ws1.Cells.Item(36, 5).Value = "LP1912"

$ws1.Cells.Item(37, 1).Value = "16:52:37"
$ws1.Cells.Item(37, 2).Value = "17:51"
$ws1.Cells.Item(37, 3).Value = "215_EL PELIGRO"
$ws1.Cells.Item(37, 4).Value = 59
$ws1.Cells.Item(37, 5).Value = "LP1912"

$ws1.Cells.Item(38, 1).Value = "16:46:42"
$ws1.Cells.Item(38, 2).Value = "18:02"
$ws1.Cells.Item(38, 3).Value = "17_ROMERO"
$ws1.Cells.Item(38, 4).Value = 76
$ws1.Cells.Item(38, 5).Value = "LP1912"

$ws1.Cells.Item(39, 1).Value = "16:52:37"
$ws1.Cells.Item(39, 2).Value = "18:03"
$ws1.Cells.Item(39, 3).Value = "17_ROMERO"
$ws1.Cells.Item(39, 4).Value = 71
$ws1.Cells.Item(39, 5).Value = "LP1912"

$ws1.Cells.Item(40, 1).Value = "16:46:42"
$ws1.Cells.Item(40, 2).Value = "18:04"
$ws1.Cells.Item(40, 3).Value = "14_ABASTO"
$ws1.Cells.Item(40, 4).Value = 78
$ws1.Cells.Item(40, 5).Value = "LP1912"

$ws1.Cells.Item(41, 1).Value = "16:52:37"
$ws1.Cells.Item(41, 2).Value = "18:14"
$ws1.Cells.Item(41, 3).Value = "10_OLMOS"
$ws1.Cells.Item(41, 4).Value = 82
$ws1.Cells.Item(41, 5).Value = "LP1912"

$ws1.Cells.Item(42, 1).Value = "16:46:42"
$ws1.Cells.Item(42, 2).Value = "18:24"
$ws1.Cells.Item(42, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(42, 4).Value = 98
$ws1.Cells.Item(42, 5).Value = "LP1912"

$ws1.Cells.Item(43, 1).Value = "17:13:30"
$ws1.Cells.Item(43, 2).Value = "18:27"
$ws1.Cells.Item(43, 3).Value = "15_ABASTO"
$ws1.Cells.Item(43, 4).Value = 74
$ws1.Cells.Item(43, 5).Value = "LP1912"

$ws1.Cells.Item(44, 1).Value = "16:46:42"
$ws1.Cells.Item(44, 2).Value = "18:34"
$ws1.Cells.Item(44, 3).Value = "14X44_ABASTO"
$ws1.Cells.Item(44, 4).Value = 108
$ws1.Cells.Item(44, 5).Value = "LP1912"

$ws1.Cells.Item(45, 1).Value = "16:46:42"
$ws1.Cells.Item(45, 2).Value = "18:38"
$ws1.Cells.Item(45, 3).Value = "17X38_ROMERO"
$ws1.Cells.Item(45, 4).Value = 112
$ws1.Cells.Item(45, 5).Value = "LP1912"

$ws1.Cells.Item(46, 1).Value = "16:46:42"
$ws1.Cells.Item(46, 2).Value = "18:41"
$ws1.Cells.Item(46, 3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(46, 4).Value = 115
$ws1.Cells.Item(46, 5).Value = "LP1912"

$ws1.Cells.Item(47, 1).Value = "17:13:30"
$ws1.Cells.Item(47, 2).Value = "18:41"
$ws1.Cells.Item(47, 3).Value = "14_ABASTO"
$ws1.Cells.Item(47, 4).Value = 88
$ws1.Cells.Item(47, 5).Value = "LP1912"

$ws1.Cells.Item(48, 1).Value = "17:13:30"
$ws1.Cells.Item(48, 2).Value = "19:01"
$ws1.Cells.Item(48, 3).Value = "17_ROMERO"
$ws1.Cells.Item(48, 4).Value = 108
$ws1.Cells.Item(48, 5).Value = "LP1912"

$ws1.Cells.Item(49, 1).Value = "17:13:30"
$ws1.Cells.Item(49, 2).Value = "19:11"
$ws1.Cells.Item(49, 3).Value = "81_EL PELIGRO"
$ws1.Cells.Item(49, 4).Value = 118
$ws1.Cells.Item(49, 5).Value = "LP1912"

# --- Sheet3: 6203-6173 ---
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 17:13:30"
$ws3.Range("A3").Value = "Total filas: 6"

$ws3.Cells.Item(10, 1).Value = "17:13:30"
$ws3.Cells.Item(10, 2).Value = "18:22"
$ws3.Cells.Item(10, 3).Value = "215C_LA PLATA"
$ws3.Cells.Item(10, 4).Value = 69
$ws3.Cells.Item(10, 5).Value = "L6203"

$ws3.Cells.Item(11, 1).Value = "17:13:30"
$ws3.Cells.Item(11, 2).Value = "19:11"
$ws3.Cells.Item(11, 3).Value = "215B_LP-P MOR-1 Y 57"
$ws3.Cells.Item(11, 4).Value = 118
$ws3.Cells.Item(11, 5).Value = "L6173"
